$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells in columns D and E keep their exact string
# representation (some look numeric/date-like, e.g. "1.000", "0.9998").
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.784.44"
$ws.Range("D3").Value = "1.955.03"
$ws.Range("E3").Value = "  +3.81%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "251.63"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").Value = "0.5964"
$ws.Range("E6").Value = "  +27.28%  "
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.3138"
$ws.Range("E8").Value = "  +8.48%  "
$ws.Range("D9").Value = "24.55"
$ws.Range("E9").Value = "  +10.21%  "
$ws.Range("D10").Value = "0.06851"
$ws.Range("E10").Value = "  +5.70%  "
$ws.Range("D11").Value = "0.8085"
$ws.Range("E11").Value = "  +11.05%  "
$ws.Range("D12").Value = "101.27"
$ws.Range("E12").Value = "  +6.16%  "
$ws.Range("D13").Value = "0.07969"
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("D14").Value = "1.937.43"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "5.342"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("D16").Value = "281.85"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").Value = "30.799.24"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "13.72"
$ws.Range("E18").Value = "  +5.52%  "
$ws.Range("D19").Value = "0.000007682"
$ws.Range("E19").Value = "  +2.79%  "
$ws.Range("D20").Value = "5.604"
$ws.Range("E20").Value = "  +6.70%  "
$ws.Range("D21").Value = "2.185.63"
$ws.Range("E21").Value = "  +2.96%  "
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "6.635"
$ws.Range("E24").Value = "  +6.17%  "
$ws.Range("D25").Value = "9.476"
$ws.Range("E25").Value = "  +4.39%  "
$ws.Range("D26").Value = "165.29"
$ws.Range("E26").Value = "  +0.96%  "
$ws.Range("D27").Value = "19.54"
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").Value = "2.083"
$ws.Range("E28").Value = "  +10.19%  "
$ws.Range("D29").Value = "0.1107"
$ws.Range("E29").Value = "  +14.17%  "
$ws.Range("D30").Value = "1.355"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "1.548"
$ws.Range("E31").Value = "  +5.42%  "
$ws.Range("D32").Value = "4.473"
$ws.Range("E32").Value = "  +4.64%  "
$ws.Range("D33").Value = "4.357"
$ws.Range("E33").Value = "  +5.64%  "
$ws.Range("D34").Value = "0.04991"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "1.178"
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("D36").Value = "0.7224"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "2.722"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  +3.93%  "
$ws.Range("D39").Value = "2.923"
$ws.Range("E39").Value = "  +3.49%  "
$ws.Range("D40").Value = "77.72"
$ws.Range("E40").Value = "  +2.92%  "
$ws.Range("D41").Value = "6.432"
$ws.Range("E41").Value = "  +4.18%  "
$ws.Range("D42").Value = "0.4505"
$ws.Range("E42").Value = "  +6.03%  "
$ws.Range("D44").Value = "0.8448"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "10.11"
$ws.Range("E46").Value = "  +6.60%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "102.69"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("D48").Value = "7.293"
$ws.Range("E48").Value = "  +4.80%  "
$ws.Range("D49").Value = "36.02"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("D50").Value = "0.4135"
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").Value = "913.77"
$ws.Range("E51").Value = "  -0.16%  "
